# Updates cryptos list data (prices / volume-1h deltas, plus a few re-ordered
# coin rows) to match the source refresh performed by the scheduled GitHub
# Actions job on 2023-05-08.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.245.85"
$ws.Range("E2").Value = "  -2.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.03"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.36"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4384"
$ws.Range("E7").Value = "  -4.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3714"
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07538"
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9412"
$ws.Range("E10").Value = "  -3.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.39"
$ws.Range("E11").Value = "  -2.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.853.98"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.742"
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.468"
$ws.Range("E14").Value = "  -2.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06866"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "82.38"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009120"
$ws.Range("E18").Value = "  -3.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.0000"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  -3.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.219.90"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.147"
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.79"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.076.71"
$ws.Range("E24").Value = "  -2.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.82"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.47"
$ws.Range("E27").Value = "  -2.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.348"
$ws.Range("E28").Value = "  -4.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.08"
$ws.Range("E29").Value = "  -2.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.737"
$ws.Range("E30").Value = "  -5.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09046"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8048"
$ws.Range("E32").Value = "  -6.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.863"
$ws.Range("E33").Value = "  -4.13%  "
$ws.Range("E34").Value = "  -5.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.964"
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.122"
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05471"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01959"
$ws.Range("E39").Value = "  -3.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.975"
$ws.Range("E40").Value = "  +8.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.147"
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5258"
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1680"
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.756"
$ws.Range("E44").Value = "  -5.32%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06779"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.049"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4884"
$ws.Range("E47").Value = "  -4.95%  "
$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000002544"
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.55"
$ws.Range("E49").Value = "  -5.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "107.83"
$ws.Range("E50").Value = "  -1.85%  "
